$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@(48, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Akşam Lisesi", "Binbaşı (Kara Kuvvetleri)", "lise")
    ,@(49, "d", 22, "Kadın", "fen-matematik", "fen-matematik", "Anadolu Lisesi (Yabancı Dille Öğretim Yapan R", "Genel Müdür Yardımcısı-Eğitim, Sağlık, Spor (", "lise")
    ,@(50, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Anadolu Lisesi (Yabancı Dille Öğretim Yapan R", "Muhtar-Köy", "lise")
    ,@(51, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Anadolu Lisesi (Yabancı Dille Öğretim Yapan R", "Binbaşı (Jandarma)", "lise")
    ,@(52, "gathaus", 22, "Kadın", "fen-matematik", "fen-matematik", "Akşam Lisesi", "Albay (Sahil Güvenlik)", "lise")
    ,@(53, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Akşam Lisesi", "Genel Müdür-Eğitim, Sağlık, Spor (Özel Sektör", "lise")
    ,@(54, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Akşam Lisesi", "(Lisans Programı 4-6 Yıllık)ıştay Başkanı", "lise")
    ,@(55, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Akşam Lisesi", "Yönetici-Elektrik, Havagazı, Su Ve Sıhhi Tesi", "lise")
    ,@(56, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Anadolu Lisesi (Yabancı Dille Öğretim Yapan R", "Ayakkabı ve Saraciye Teknolojisi Öğretmeni", "lise")
    ,@(57, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Akşam Lisesi", "Genel Müdür-İnşaat Ve İmalat Sanayi (Özel Sek", "lise")
    ,@(58, "gathaus", 22, "Kadın", "fen-matematik", "fen-matematik", "Anadolu Lisesi (Yabancı Dille Öğretim Yapan R", "Genel Müdür-Bankalar Ve Mali Müesseseler (Öze", "lise")
    ,@(59, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Anadolu Lisesi (Yabancı Dille Öğretim Yapan R", "Halk Oyunları Oyuncusu (Siirt Yöresi)", "lise")
    ,@(60, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Anadolu Lisesi (Yabancı Dille Öğretim Yapan R", "Genel Müdür-Bankalar Ve Mali Müesseseler (Öze", "lise")
    ,@(61, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Akşam Lisesi", "Yönetici-Elektrik, Havagazı, Su Ve Sıhhi Tesi", "lise")
    ,@(62, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Anadolu Lisesi (Yabancı Dille Öğretim Yapan R", "Genel Müdür-İnşaat Ve İmalat Sanayi (Özel Sek", "lise")
    ,@(63, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Anadolu Lisesi (Yabancı Dille Öğretim Yapan R", "Genel Müdür Yardımcısı-Bankalar Ve Mali Müess", "lise")
    ,@(64, "gathaus", 22, "Kadın", "fen-matematik", "fen-matematik", "Akşam Lisesi", "Genel Müdür Yardımcısı-İnşaat Ve İmalat Sanay", "lise")
    ,@(65, "Rıza Mert Yağcı", 22, "Kadın", "fen-matematik", "fen-matematik", "Anadolu Lisesi (Yabancı Dille Öğretim Yapan R", "Genel Müdür-İnşaat Ve İmalat Sanayi (Özel Sek", "lise")
)

$startRow = 48
for ($idx = 0; $idx -lt $rows.Count; $idx++) {
    $r = $startRow + $idx
    $data = $rows[$idx]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
}

Write-Output "Done adding rows 48-65"
